$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.499.57'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.60'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.84'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6295'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07533'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2977'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.39'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07724'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.856.91'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.002'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.57'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009778'
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.112.89'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.247'
$ws.Range('E18').Value = '  +2.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.564.59'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '233.39'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.49'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.654'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.76'
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1389'
$ws.Range('E26').Value = '  -2.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.459'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.70'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.477'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05944'
$ws.Range('E30').Value = '  -3.80%  '
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.112'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.036'
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.890'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7194'
$ws.Range('E36').Value = '  -1.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.592'
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.797'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.239.37'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01796'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9084'
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.142'
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.053.43'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.40'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.20'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.397'
$ws.Range('E47').Value = '  +9.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000119'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.177'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4046'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.701'
$ws.Range('E51').Value = '  +2.06%  '
